# Lastenheft.docx edit script
# Applies the textual corrections described by the commit diff:
#  - "Trafic-Nois-Detector" -> "Traffic-Noise-Detector" (typo fixes in the title)
#  - "50Km/h" -> "50 Km/h" (add missing space)
#  - drop ", für den Privatgebrauch," aside
#  - "an einen PC übersandt" -> "an Peripheriegerät gesendet"
#  - "Dieser kann sich" -> "Der Benutzer kann sich"
#  - "48kHz bei 24Bit" -> "48 kHz bei 24 Bit" (add missing spaces)
#  - re-insert the "_GoBack" bookmark after "... Option auf zwei weitere"

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $rng.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Title: "Trafic-Nois-Detector" -> "Traffic-Noise-Detector"
Replace-Text "Trafic-Nois-Detector" "Traffic-Noise-Detector"

# 2. "innerorts bis 50Km/h" -> "innerorts bis 50 Km/h"
Replace-Text "innerorts bis 50Km/h" "innerorts bis 50 Km/h"

# 3. Drop the "für den Privatgebrauch" aside after "Rucksack"
Replace-Text "Rucksack, für den Privatgebrauch, verstaubar" "Rucksack verstaubar"

# 4. "an einen PC übersandt," -> "an Peripheriegerät gesendet,"
Replace-Text "einen PC übersandt," "Peripheriegerät gesendet,"

# 5. "Dieser kann sich in einem" -> "Der Benutzer kann sich in einem"
Replace-Text "Dieser kann sich in einem" "Der Benutzer kann sich in einem"

# 6. "48kHz bei 24Bit durchgeführt" -> "48 kHz bei 24 Bit durchgeführt"
Replace-Text "Messfrequenz von 48kHz bei 24Bit durchgeführt" "Messfrequenz von 48 kHz bei 24 Bit durchgeführt"

# 7. Re-add the "_GoBack" bookmark right after "... Option auf zwei weitere"
$bm = $d.Content
$bm.Find.Execute("Option auf zwei weitere", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bm.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bm)
